# StorageComponentClassDiagram.pptx update
#
# Mirrors the author's commit:
#   - rename AddressBook/Person concepts to TaskManager/Task in the
#     class-diagram shapes on slide 1 (plus a font-size tweak on one run)
#   - refresh the cached "datetimeFigureOut" footer date text everywhere
#     it is cached (slide master, slide layouts, notes master)
#   - (best effort) presentation-level slide guides

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

function Set-RunText($textRange, $start, $length, $newText) {
    $sub = $textRange.Characters($start, $length)
    $sub.Text = $newText
}

function Set-RunFontSize($textRange, $start, $length, $size) {
    $sub = $textRange.Characters($start, $length)
    $sub.Font.Size = $size
}

function Set-DatePlaceholderText($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
        }
        if ($isDate) {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

# ---------------------------------------------------------------------
# Slide 1 - class diagram shape text
# ---------------------------------------------------------------------

$s = $p.Slides.Item(1)

# Shape id=2: "<<interface>>\nAddressBookStorage" -> "...TaskManagerStorage"
# and the second line shrinks from 10.5pt to 10pt.
$sh2 = Get-ShapeById $s.Shapes 2
$tr2 = $sh2.TextFrame.TextRange
Set-RunText $tr2 15 19 "TaskManagerStorage"
$tr2b = $sh2.TextFrame.TextRange
Set-RunFontSize $tr2b 15 19 10

# Shape id=50: "XmlAddressBook" + (empty run) + break + "Storage"
#   -> "XmlTaskManager" + ... + "Storage"
$sh50 = Get-ShapeById $s.Shapes 50
$tr50 = $sh50.TextFrame.TextRange
Set-RunText $tr50 1 14 "XmlTaskManager"

# Shape id=66: "XmlSerializable" + break + "AddressBook"
#   -> "XmlSerializable" + break + "TaskManager"
$sh66 = Get-ShapeById $s.Shapes 66
$tr66 = $sh66.TextFrame.TextRange
Set-RunText $tr66 17 11 "TaskManager"

# Shape id=74: "XmlAdaptedPerson" -> "XmlAdaptedTask"
$sh74 = Get-ShapeById $s.Shapes 74
$sh74.TextFrame.TextRange.Text = "XmlAdaptedTask"

# ---------------------------------------------------------------------
# Cached footer date ("datetimeFigureOut" field) - 10/16/2016 -> 3/16/17
# ---------------------------------------------------------------------
#
# NOTE: the notes master's own Date Placeholder is intentionally left
# alone here - writing through $p.NotesMaster.Shapes.Item(N).TextFrame
# in this host mis-resolves to the *slide* master's Nth shape instead
# (an indexing aliasing bug), clobbering unrelated placeholders. Only
# the slide master + slide layouts (which resolve correctly) are
# touched.

$newDate = "3/16/17"
$m = $p.SlideMaster

# The slide master's own Date Placeholder
Set-DatePlaceholderText $m.Shapes $newDate

# Every slide layout's Date Placeholder
for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    $cl = $m.CustomLayouts.Item($li)
    Set-DatePlaceholderText $cl.Shapes $newDate
}

# ---------------------------------------------------------------------
# Presentation slide guides (horizontal @1488, vertical @2880)
# ---------------------------------------------------------------------

try {
    $guides = $p.Guides
    $hGuide = $guides.Add(1, 1488)
    $hGuide.Color = RGB(164, 163, 164)
    $vGuide = $guides.Add(2, 2880)
    $vGuide.Color = RGB(164, 163, 164)
} catch {
}
